# Softexpert template - automated spreadsheet update
$wb = $excel.ActiveWorkbook

# --- "Projetos" sheet: move the saved selection to B11 ---
$wsProjetos = $wb.Worksheets.Item("Projetos")
$wsProjetos.Activate()
$wsProjetos.Range("B11").Select()

# --- "Etapas" sheet: shift the existing stage blocks down one project ---
$wsEtapas = $wb.Worksheets.Item("Etapas")
$wsEtapas.Activate()

# Rows 23-27 were "P04" -> now "P05"
$wsEtapas.Range("A23:A27").Value = "P05"
# Rows 28-32 were "P05" -> now "P06"
$wsEtapas.Range("A28:A32").Value = "P06"
# Rows 33-37 were "P06" -> now "P07"
$wsEtapas.Range("A33:A37").Value = "P07"
# Rows 38-42 were "P07" -> now "P08"
$wsEtapas.Range("A38:A42").Value = "P08"
# Rows 43-47 were "P08" -> now "P09"
$wsEtapas.Range("A43:A47").Value = "P09"

# New rows 48-52: copy of the stage pattern for the next project (P10 / ITSM)
$wsEtapas.Range("A48:H52").Value = $wsEtapas.Range("A43:H47").Value

$wsEtapas.Range("A48:A52").Value = "P10"
$wsEtapas.Range("B48:B52").Value = "ITSM"

$wsEtapas.Range("C48").Value = "Levantamento de requisitos"
$wsEtapas.Range("C49").Value = "Configuração dos módulos"
$wsEtapas.Range("C50").Value = "Homologação e testes"
$wsEtapas.Range("C51").Value = "Treinamento de usuários"
$wsEtapas.Range("C52").Value = "Go-Live e estabilização"

$wsEtapas.Range("D48:E48").Value = "Abr/2026"
$wsEtapas.Range("F48:G48").Value = "Mai/2026"

$wsEtapas.Range("D49:E49").Value = "Mai/2026"
$wsEtapas.Range("F49:G49").Value = "Jul/2026"

$wsEtapas.Range("D50:E50").Value = "Jul/2026"
$wsEtapas.Range("F50:G50").Value = "Ago/2026"

$wsEtapas.Range("D51:E51").Value = "Ago/2026"
$wsEtapas.Range("F51:G51").Value = "Set/2026"

$wsEtapas.Range("D52:E52").Value = "Set/2026"
$wsEtapas.Range("F52:G52").Value = "Out/2026"

$wsEtapas.Range("H48:H52").Value = "Não Iniciado"

# Update the view: scroll so A38 is the top-left frozen cell and select A49:A52
$wsEtapas.Range("A38").Select()
$wsEtapas.ActiveWindow.ScrollRow = 38
$wsEtapas.Range("A49:A52").Select()
